$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '63.885.32'
Set-TextValue $ws.Range("E2") '  +0.03%  '
Set-TextValue $ws.Range("D3") '3.134.22'
Set-TextValue $ws.Range("E3") '  +0.50%  '
Set-TextValue $ws.Range("E4") '  +0.06%  '
Set-TextValue $ws.Range("D5") '589.26'
Set-TextValue $ws.Range("E5") '  +0.56%  '
Set-TextValue $ws.Range("D6") '145.00'
Set-TextValue $ws.Range("E6") '  -0.82%  '
Set-TextValue $ws.Range("D8") '3.127.08'
Set-TextValue $ws.Range("E8") '  +0.50%  '
Set-TextValue $ws.Range("E9") '  -0.31%  '
Set-TextValue $ws.Range("E10") '  -0.86%  '
Set-TextValue $ws.Range("E11") '  +2.39%  '
Set-TextValue $ws.Range("D12") '0.457'
Set-TextValue $ws.Range("E12") '  -1.82%  '
Set-TextValue $ws.Range("E13") '  -2.53%  '
Set-TextValue $ws.Range("D14") '37.28'
Set-TextValue $ws.Range("E14") '  +0.79%  '
Set-TextValue $ws.Range("D15") '3.656.77'
Set-TextValue $ws.Range("E15") '  +0.61%  '
Set-TextValue $ws.Range("E16") '  -1.26%  '
Set-TextValue $ws.Range("E17") '  +2.36%  '
Set-TextValue $ws.Range("D18") '63.736.02'
Set-TextValue $ws.Range("E18") '  -0.02%  '
Set-TextValue $ws.Range("D19") '3.135.19'
Set-TextValue $ws.Range("E19") '  +0.73%  '
Set-TextValue $ws.Range("D20") '466.92'
Set-TextValue $ws.Range("E20") '  +0.55%  '
Set-TextValue $ws.Range("D21") '14.32'
Set-TextValue $ws.Range("E21") '  +0.12%  '
Set-TextValue $ws.Range("D22") '0.730'
Set-TextValue $ws.Range("E22") '  -0.05%  '
Set-TextValue $ws.Range("D23") '7.53'
Set-TextValue $ws.Range("E23") '  -0.02%  '
Set-TextValue $ws.Range("D24") '81.55'
Set-TextValue $ws.Range("E24") '  -0.64%  '
Set-TextValue $ws.Range("D25") '12.94'
Set-TextValue $ws.Range("E25") '  -1.44%  '
Set-TextValue $ws.Range("E26") '  +7.14%  '
Set-TextValue $ws.Range("E27") '  +0.08%  '
Set-TextValue $ws.Range("D28") '9.85'
Set-TextValue $ws.Range("E28") '  +10.02%  '
Set-TextValue $ws.Range("E29") '  +8.59%  '
Set-TextValue $ws.Range("D30") '2.70'
Set-TextValue $ws.Range("E30") '  +0.27%  '
Set-TextValue $ws.Range("D31") '2.23'
Set-TextValue $ws.Range("E31") '  +0.14%  '
Set-TextValue $ws.Range("E32") '  +0.13%  '
Set-TextValue $ws.Range("D33") '27.64'
Set-TextValue $ws.Range("E33") '  +2.41%  '
Set-TextValue $ws.Range("E34") '  +0.16%  '
Set-TextValue $ws.Range("D35") '0.0₃0849'
Set-TextValue $ws.Range("E35") '  -3.11%  '
Set-TextValue $ws.Range("E36") '  +1.19%  '
Set-TextValue $ws.Range("D37") '6.14'
Set-TextValue $ws.Range("E37") '  +0.88%  '
Set-TextValue $ws.Range("E39") '  -6.14%  '
Set-TextValue $ws.Range("D40") '51.30'
Set-TextValue $ws.Range("E40") '  +0.69%  '
Set-TextValue $ws.Range("D41") '9.33'
Set-TextValue $ws.Range("E41") '  +7.38%  '
Set-TextValue $ws.Range("D42") '452.85'
Set-TextValue $ws.Range("E42") '  +0.69%  '
Set-TextValue $ws.Range("E43") '  +5.54%  '
Set-TextValue $ws.Range("E44") '  -0.05%  '
Set-TextValue $ws.Range("D45") '2.915.34'
Set-TextValue $ws.Range("E45") '  +1.01%  '
Set-TextValue $ws.Range("D46") '40.13'
Set-TextValue $ws.Range("E46") '  +12.09%  '
Set-TextValue $ws.Range("D47") '0.108'
Set-TextValue $ws.Range("E47") '  -2.95%  '
Set-TextValue $ws.Range("D48") '132.76'
Set-TextValue $ws.Range("E48") '  +6.70%  '
# Row 50/51 swap: ThetaToken/Stellar positions exchanged
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D50") '0.111'
Set-TextValue $ws.Range("E50") '  -0.69%  '

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D51") '2.23'
Set-TextValue $ws.Range("E51") '  +2.30%  '
